$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Rows 15 and 16 swap content: Litecoin (was row15) <-> ShibaInu (was row16)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D15") "0.00001000"
$ws.Range("E15").Value = "  +15.74%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D16") "83.09"
$ws.Range("E16").Value = "  +1.32%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$ws.Range("D2").Value = "29.100.02"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.837.07"
$ws.Range("E3").Value = "  +0.44%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue $ws.Range("D5") "243.47"
$ws.Range("E5").Value = "  +0.82%  "
Set-TextValue $ws.Range("D6") "0.6275"
$ws.Range("E6").Value = "  -0.65%  "
Set-TextValue $ws.Range("D7") "1.002"
$ws.Range("E7").Value = "  +0.18%  "
Set-TextValue $ws.Range("D8") "0.07593"
$ws.Range("E8").Value = "  +3.53%  "
Set-TextValue $ws.Range("D9") "0.2929"
$ws.Range("E9").Value = "  -0.05%  "
Set-TextValue $ws.Range("D10") "22.63"
$ws.Range("E10").Value = "  -1.09%  "
Set-TextValue $ws.Range("D11") "0.07753"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "1.846.63"
$ws.Range("E12").Value = "  +0.94%  "
Set-TextValue $ws.Range("D13") "4.975"
$ws.Range("E13").Value = "  -0.26%  "
Set-TextValue $ws.Range("D14") "0.6656"
$ws.Range("E14").Value = "  +0.40%  "
Set-TextValue $ws.Range("D17") "6.078"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "29.117.43"
$ws.Range("E18").Value = "  +0.45%  "
Set-TextValue $ws.Range("D19") "227.15"
$ws.Range("E19").Value = "  +1.22%  "
Set-TextValue $ws.Range("D20") "12.42"
$ws.Range("E20").Value = "  +0.15%  "
Set-TextValue $ws.Range("D21") "1.003"
$ws.Range("E21").Value = "  +0.29%  "
Set-TextValue $ws.Range("D22") "7.223"
$ws.Range("E22").Value = "  +1.18%  "
Set-TextValue $ws.Range("D23") "1.002"
$ws.Range("E23").Value = "  +0.15%  "
Set-TextValue $ws.Range("D24") "159.66"
$ws.Range("E24").Value = "  +1.17%  "
Set-TextValue $ws.Range("D25") "8.529"
$ws.Range("E25").Value = "  +0.90%  "
Set-TextValue $ws.Range("D26") "0.1386"
$ws.Range("E26").Value = "  +1.31%  "
Set-TextValue $ws.Range("D27") "17.96"
$ws.Range("E27").Value = "  +0.39%  "
Set-TextValue $ws.Range("D28") "1.495"
$ws.Range("E28").Value = "  -0.71%  "
Set-TextValue $ws.Range("D29") "4.107"
$ws.Range("E29").Value = "  +0.42%  "
Set-TextValue $ws.Range("D30") "4.025"
$ws.Range("E30").Value = "  +0.07%  "
Set-TextValue $ws.Range("D31") "1.196"
$ws.Range("E31").Value = "  -0.52%  "
Set-TextValue $ws.Range("D32") "0.05262"
$ws.Range("E32").Value = "  -0.71%  "
Set-TextValue $ws.Range("D33") "1.845"
$ws.Range("E33").Value = "  +1.11%  "
Set-TextValue $ws.Range("D34") "0.7356"
$ws.Range("E34").Value = "  -0.50%  "
Set-TextValue $ws.Range("D35") "1.137"
$ws.Range("E35").Value = "  -1.23%  "
Set-TextValue $ws.Range("D36") "2.705"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").Value = "1.240.58"
$ws.Range("E37").Value = "  -4.10%  "
Set-TextValue $ws.Range("D38") "2.766"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  +1.03%  "
Set-TextValue $ws.Range("D41") "0.8982"
$ws.Range("E41").Value = "  +0.73%  "
Set-TextValue $ws.Range("D42") "1.002"
$ws.Range("E42").Value = "  +0.22%  "
Set-TextValue $ws.Range("D43") "102.22"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "1.986.46"
$ws.Range("E44").Value = "  +0.46%  "
Set-TextValue $ws.Range("D45") "0.00000000124"
Set-TextValue $ws.Range("D46") "64.47"
$ws.Range("E46").Value = "  +0.52%  "
Set-TextValue $ws.Range("D47") "0.5120"
$ws.Range("E47").Value = "  -0.37%  "
Set-TextValue $ws.Range("D48") "0.4044"
$ws.Range("E48").Value = "  +1.46%  "
Set-TextValue $ws.Range("D49") "8.866"
$ws.Range("E49").Value = "  +1.99%  "
Set-TextValue $ws.Range("D50") "0.05759"
Set-TextValue $ws.Range("D51") "6.696"
$ws.Range("E51").Value = "  +0.23%  "
